# Daily attendance processing - reorder "Recorded By" contributor lists.
#
# For each data row in the "Recorded By" column (G), the comma-separated
# list of contributors is rotated so that the last-listed contributor
# moves to the front (e.g. "System, dnasr281@gmail.com" becomes
# "dnasr281@gmail.com, System"). Single-value cells are left untouched,
# and the specific combination "System, admin@admin.com" is preserved
# as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "" -and $val -ne "System, admin@admin.com") {
        $parts = $val.Split(",")
        foreach ($i in 0..($parts.Count - 1)) {
            $parts[$i] = $parts[$i].Trim()
        }

        if ($parts.Count -gt 1) {
            $n = $parts.Count
            $last = $parts[$n - 1]
            $rest = $parts[0..($n - 2)]
            $reordered = @($last) + $rest
            $newVal = [string]::Join(", ", $reordered)

            if ($newVal -ne $val) {
                $cell.Value2 = $newVal
            }
        }
    }
}
